# edit.ps1 - apply "updated log book pages" changes via Word COM interop
#
# Changes applied:
#  1. "SEM III & IV :" -> "SEM V & VI :" (with proofErr gramStart/gramEnd wrapping "VI :")
#  2. "Academic Year: 20__ to 20__ " (tab-underlines) -> "Academic Year: 2021 to 2022" (typed digits)
#  3. "Date:" line gets an extra duplicated spacing run (a stray space run) inserted after "Date:"
#  4. A new "Team Member 4:" paragraph + a new blank paragraph are inserted after "Team Member 3:"

$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1 & 2. Rebuild the "SEM ... WEEK ... Academic Year ..." paragraph (paragraph 6)
# ---------------------------------------------------------------------------
$pkgHeader = '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>'
$pkgFooter = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$semWeekXml = $pkgHeader + '<w:p><w:pPr><w:pStyle w:val="BodyText"/><w:tabs><w:tab w:val="left" w:pos="3123"/><w:tab w:val="left" w:pos="6497"/><w:tab w:val="left" w:pos="8626"/><w:tab w:val="left" w:pos="9412"/></w:tabs><w:spacing w:before="90"/><w:ind w:left="115"/></w:pPr><w:r><w:t>SEM V</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>&amp;</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t xml:space="preserve">VI </w:t></w:r><w:r><w:t>:</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:spacing w:val="-19"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>WEEK</w:t></w:r><w:r><w:rPr><w:spacing w:val="2"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>–</w:t></w:r><w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:u w:val="single"/></w:rPr><w:tab/></w:r><w:r><w:tab/><w:t>Academic</w:t></w:r><w:r><w:rPr><w:spacing w:val="-9"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:spacing w:val="-6"/></w:rPr><w:t>Year:</w:t></w:r><w:r><w:rPr><w:spacing w:val="1"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>20</w:t></w:r><w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t>21</w:t></w:r><w:r><w:t>to</w:t></w:r><w:r><w:rPr><w:spacing w:val="-1"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>20</w:t></w:r><w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t>22</w:t></w:r><w:r><w:rPr><w:u w:val="single"/></w:rPr><w:tab/></w:r></w:p>' + $pkgFooter

$semWeekPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs.Item($i)
    if ($candidate.Range.Text -like "SEM III*") {
        $semWeekPara = $candidate
        break
    }
}
$semWeekPara.Range.InsertXML($semWeekXml)

# ---------------------------------------------------------------------------
# 3. "Date:" paragraph - insert a duplicated spacing run right after "Date:"
# ---------------------------------------------------------------------------
$datePara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs.Item($i)
    if ($candidate.Range.Text -like "Date:*From*To:*") {
        $datePara = $candidate
        break
    }
}
$dateXml = $pkgHeader + '<w:p><w:pPr><w:pStyle w:val="BodyText"/><w:tabs><w:tab w:val="left" w:pos="2773"/><w:tab w:val="left" w:pos="4766"/></w:tabs><w:spacing w:before="90"/><w:ind w:left="115"/></w:pPr><w:r><w:t>Date:</w:t></w:r><w:r><w:rPr><w:spacing w:val="-1"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:spacing w:val="-1"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>From</w:t></w:r><w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:u w:val="single"/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:spacing w:val="-6"/></w:rPr><w:t>To:</w:t></w:r><w:r><w:rPr><w:spacing w:val="1"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:u w:val="single"/></w:rPr><w:tab/></w:r></w:p>' + $pkgFooter
$datePara.Range.InsertXML($dateXml)

# ---------------------------------------------------------------------------
# 4. Insert "Team Member 4:" paragraph + a blank paragraph after "Team Member 3:"
# ---------------------------------------------------------------------------
$teamMember3 = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs.Item($i)
    if ($candidate.Range.Text -like "Team Member 3:*") {
        $teamMember3 = $candidate
        break
    }
}
$tm3Index = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "Team Member 3:*") {
        $tm3Index = $i
        break
    }
}

$teamMember3.Range.InsertParagraphAfter()

$teamMember4Para = $d.Paragraphs.Item($tm3Index + 1)
$teamMember4Xml = $pkgHeader + '<w:p><w:pPr><w:pStyle w:val="BodyText"/><w:ind w:left="115"/></w:pPr><w:r><w:rPr><w:spacing w:val="-5"/></w:rPr><w:t xml:space="preserve">Team </w:t></w:r><w:r><w:t>Member</w:t></w:r><w:r><w:rPr><w:spacing w:val="3"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>4</w:t></w:r><w:r><w:t>:</w:t></w:r></w:p>' + $pkgFooter
$teamMember4Para.Range.InsertXML($teamMember4Xml)

$teamMember4Para = $d.Paragraphs.Item($tm3Index + 1)
$teamMember4Para.Range.InsertParagraphAfter()

$blankPara = $d.Paragraphs.Item($tm3Index + 2)
$blankXml = $pkgHeader + '<w:p><w:pPr><w:pStyle w:val="BodyText"/><w:ind w:left="115"/></w:pPr></w:p>' + $pkgFooter
$blankPara.Range.InsertXML($blankXml)

Write-Output "edit complete; paragraph count = $($d.Paragraphs.Count)"
